# Atualiza dados de faturamento da Bibi (coluna L = dia 11) e os totais (coluna AG)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linha 2 - Bibi Cell Mundi
$ws.Range("L2").Value = 20151.8
$ws.Range("AG2").Value = 106783.53

# Linha 3 - Bibi Cell Vieiralves
$ws.Range("L3").Value = 4295
$ws.Range("AG3").Value = 56680.2

# Linha 4 - Bibi Cell Manauara
$ws.Range("L4").Value = 2671
$ws.Range("AG4").Value = 33510.15

# Linha 5 - Bibi Cell Ponta Negra
$ws.Range("L5").Value = 5724
$ws.Range("AG5").Value = 31041.33

# Linha 6 - total
$ws.Range("L6").Value = 32841.8
$ws.Range("AG6").Value = 228015.21
